$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.203781008720398
$ws.Range("B1").Value = 1.621781587600708
$ws.Range("C1").Value = 6.941516399383545
$ws.Range("D1").Value = 2.213505983352661
$ws.Range("E1").Value = 1.177436232566833
